$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a few "Razon social" / "Nombre Fantasia" entries where a comma was
# mis-scraped as part of the company/person name; normalize the stray comma
# to a period (matches the corrected source data). ---
$nameFixes = @(
    @{ Cell = "E160"; Value = "GIMENEZ ANIBAL. FALISTOCCO MARISA DANIELA SH" },
    @{ Cell = "E211"; Value = "PARPAGNOLI. PEDRO RICARDO" },
    @{ Cell = "F211"; Value = "PARPAGNOLI. PEDRO RICARDO" },
    @{ Cell = "E214"; Value = "RICCOTTI. MARIANA EDITH" }
)
foreach ($fix in $nameFixes) {
    $ws.Range($fix.Cell).Value2 = $fix.Value
}

# --- Fix "Importe" (column H) formatting: values were scraped with
# Spanish/Argentine number formatting (thousands separator "." and decimal
# comma ","), e.g. "18.300,00". Convert them to plain decimal-point numeric
# text, e.g. "18300.00", while keeping them as text cells. ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Value2
    if ($null -eq $old -or $old -eq "") { continue }
    $old = [string]$old
    $new = $old.Replace(".", "").Replace(",", ".")
    if ($new -ne $old) {
        # Leading apostrophe forces the numeric-looking text to stay text
        # instead of being re-interpreted as a number.
        $cell.Value2 = "'" + $new
    }
}
